$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 15 and row 16 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(15, 2).Value = 7004588
$ws.Cells.Item(16, 2).Value = 7004589
$ws.Cells.Item(15, 3).Value = "Qatar Stars League"
$ws.Cells.Item(16, 3).Value = "Qatar Stars League"
$ws.Cells.Item(15, 4).Value = "Qatar Stars League"
$ws.Cells.Item(16, 4).Value = "Qatar Stars League"
$ws.Cells.Item(15, 5).Value = 45170.59375
$ws.Cells.Item(16, 5).Value = 45170.59375
$ws.Cells.Item(15, 6).Value = "Umm Salal"
$ws.Cells.Item(16, 6).Value = "AlMuaidar"
$ws.Cells.Item(15, 7).Value = "Qatar SC Doha"
$ws.Cells.Item(16, 7).Value = "Al Markhiya"
$ws.Cells.Item(15, 8).Value = 2
$ws.Cells.Item(16, 8).Value = 5
$ws.Cells.Item(15, 9).Value = 2
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(15, 10).Value = "D"
$ws.Cells.Item(16, 10).Value = "H"
$ws.Cells.Item(15, 11).Value = 2.9
$ws.Cells.Item(16, 11).Value = 2.4
$ws.Cells.Item(15, 12).Value = 3.3
$ws.Cells.Item(16, 12).Value = 3.6
$ws.Cells.Item(15, 13).Value = 2.25
$ws.Cells.Item(16, 13).Value = 2.375
$ws.Cells.Item(15, 14).Value = 3.75
$ws.Cells.Item(16, 14).Value = 2.5
$ws.Cells.Item(15, 15).Value = 3.4
$ws.Cells.Item(16, 15).Value = 3.5
$ws.Cells.Item(15, 16).Value = 1.909
$ws.Cells.Item(16, 16).Value = 2.3
$ws.Cells.Item(15, 17).Value = 0.5
$ws.Cells.Item(16, 17).Value = 0
$ws.Cells.Item(15, 18).Value = 1.85
$ws.Cells.Item(16, 18).Value = 1.975
$ws.Cells.Item(15, 19).Value = 1.95
$ws.Cells.Item(16, 19).Value = 1.825
$ws.Cells.Item(15, 20).Value = 2.5
$ws.Cells.Item(16, 20).Value = 3
$ws.Cells.Item(15, 21).Value = 1.85
$ws.Cells.Item(16, 21).Value = 2
$ws.Cells.Item(15, 22).Value = 1.95
$ws.Cells.Item(16, 22).Value = 1.8
$ws.Cells.Item(15, 23).Value = -1
$ws.Cells.Item(16, 23).Value = 1.5
$ws.Cells.Item(15, 24).Value = 2.4
$ws.Cells.Item(16, 24).Value = -1
$ws.Cells.Item(15, 25).Value = -1
$ws.Cells.Item(16, 25).Value = -1
$ws.Cells.Item(15, 26).Value = 0.8500000000000001
$ws.Cells.Item(16, 26).Value = 0.9750000000000001
$ws.Cells.Item(15, 27).Value = -1
$ws.Cells.Item(16, 27).Value = -1
$ws.Cells.Item(15, 28).Value = 0.8500000000000001
$ws.Cells.Item(16, 28).Value = 1
$ws.Cells.Item(15, 29).Value = -1
$ws.Cells.Item(16, 29).Value = -1

# Swap row 18 and row 19 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(18, 2).Value = 7003585
$ws.Cells.Item(19, 2).Value = 7004591
$ws.Cells.Item(18, 3).Value = "Qatar Stars League"
$ws.Cells.Item(19, 3).Value = "Qatar Stars League"
$ws.Cells.Item(18, 4).Value = "Qatar Stars League"
$ws.Cells.Item(19, 4).Value = "Qatar Stars League"
$ws.Cells.Item(18, 5).Value = 45171.59375
$ws.Cells.Item(19, 5).Value = 45171.59375
$ws.Cells.Item(18, 6).Value = "Al Sadd"
$ws.Cells.Item(19, 6).Value = "AlShamal SC"
$ws.Cells.Item(18, 7).Value = "AlWakrah SC"
$ws.Cells.Item(19, 7).Value = "AlRayyan SC"
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(19, 8).Value = 3
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(19, 9).Value = 4
$ws.Cells.Item(18, 10).Value = "D"
$ws.Cells.Item(19, 10).Value = "A"
$ws.Cells.Item(18, 11).Value = 1.615
$ws.Cells.Item(19, 11).Value = 4.5
$ws.Cells.Item(18, 12).Value = 4
$ws.Cells.Item(19, 12).Value = 4.2
$ws.Cells.Item(18, 13).Value = 4.333
$ws.Cells.Item(19, 13).Value = 1.55
$ws.Cells.Item(18, 14).Value = 1.533
$ws.Cells.Item(19, 14).Value = 3.3
$ws.Cells.Item(18, 15).Value = 4.2
$ws.Cells.Item(19, 15).Value = 3.8
$ws.Cells.Item(18, 16).Value = 5
$ws.Cells.Item(19, 16).Value = 1.85
$ws.Cells.Item(18, 17).Value = -1
$ws.Cells.Item(19, 17).Value = 0.5
$ws.Cells.Item(18, 18).Value = 1.8
$ws.Cells.Item(19, 18).Value = 1.85
$ws.Cells.Item(18, 19).Value = 2
$ws.Cells.Item(19, 19).Value = 1.95
$ws.Cells.Item(18, 20).Value = 3.5
$ws.Cells.Item(19, 20).Value = 2.75
$ws.Cells.Item(18, 21).Value = 1.925
$ws.Cells.Item(19, 21).Value = 1.85
$ws.Cells.Item(18, 22).Value = 1.875
$ws.Cells.Item(19, 22).Value = 1.95
$ws.Cells.Item(18, 23).Value = -1
$ws.Cells.Item(19, 23).Value = -1
$ws.Cells.Item(18, 24).Value = 3.2
$ws.Cells.Item(19, 24).Value = -1
$ws.Cells.Item(18, 25).Value = -1
$ws.Cells.Item(19, 25).Value = 0.8500000000000001
$ws.Cells.Item(18, 26).Value = -1
$ws.Cells.Item(19, 26).Value = -1
$ws.Cells.Item(18, 27).Value = 1
$ws.Cells.Item(19, 27).Value = 0.95
$ws.Cells.Item(18, 28).Value = -1
$ws.Cells.Item(19, 28).Value = 0.8500000000000001
$ws.Cells.Item(18, 29).Value = 0.875
$ws.Cells.Item(19, 29).Value = -1

# Swap row 24 and row 25 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(24, 2).Value = 7004592
$ws.Cells.Item(25, 2).Value = 7004593
$ws.Cells.Item(24, 3).Value = "Qatar Stars League"
$ws.Cells.Item(25, 3).Value = "Qatar Stars League"
$ws.Cells.Item(24, 4).Value = "Qatar Stars League"
$ws.Cells.Item(25, 4).Value = "Qatar Stars League"
$ws.Cells.Item(24, 5).Value = 45193.58333333334
$ws.Cells.Item(25, 5).Value = 45193.58333333334
$ws.Cells.Item(24, 6).Value = "AlMuaidar"
$ws.Cells.Item(25, 6).Value = "Al Markhiya"
$ws.Cells.Item(24, 7).Value = "AlShamal SC"
$ws.Cells.Item(25, 7).Value = "AlWakrah SC"
$ws.Cells.Item(24, 8).Value = 2
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 2
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(24, 10).Value = "D"
$ws.Cells.Item(25, 10).Value = "A"
$ws.Cells.Item(24, 11).Value = 2.4
$ws.Cells.Item(25, 11).Value = 4
$ws.Cells.Item(24, 12).Value = 3.6
$ws.Cells.Item(25, 12).Value = 4
$ws.Cells.Item(24, 13).Value = 2.4
$ws.Cells.Item(25, 13).Value = 1.727
$ws.Cells.Item(24, 14).Value = 2.4
$ws.Cells.Item(25, 14).Value = 4.75
$ws.Cells.Item(24, 15).Value = 3.6
$ws.Cells.Item(25, 15).Value = 4.333
$ws.Cells.Item(24, 16).Value = 2.45
$ws.Cells.Item(25, 16).Value = 1.571
$ws.Cells.Item(24, 17).Value = 0
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(24, 18).Value = 1.925
$ws.Cells.Item(25, 18).Value = 1.85
$ws.Cells.Item(24, 19).Value = 1.875
$ws.Cells.Item(25, 19).Value = 1.95
$ws.Cells.Item(24, 20).Value = 2.75
$ws.Cells.Item(25, 20).Value = 3
$ws.Cells.Item(24, 21).Value = 1.75
$ws.Cells.Item(25, 21).Value = 1.825
$ws.Cells.Item(24, 22).Value = 1.95
$ws.Cells.Item(25, 22).Value = 1.975
$ws.Cells.Item(24, 23).Value = -1
$ws.Cells.Item(25, 23).Value = -1
$ws.Cells.Item(24, 24).Value = 2.6
$ws.Cells.Item(25, 24).Value = -1
$ws.Cells.Item(24, 25).Value = -1
$ws.Cells.Item(25, 25).Value = 0.571
$ws.Cells.Item(24, 26).Value = 0
$ws.Cells.Item(25, 26).Value = -1
$ws.Cells.Item(24, 27).Value = -0
$ws.Cells.Item(25, 27).Value = 0.95
$ws.Cells.Item(24, 28).Value = 0.75
$ws.Cells.Item(25, 28).Value = 0
$ws.Cells.Item(24, 29).Value = -1
$ws.Cells.Item(25, 29).Value = -0

# Swap row 40 and row 41 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(40, 2).Value = 7004607
$ws.Cells.Item(41, 2).Value = 7004604
$ws.Cells.Item(40, 3).Value = "Qatar Stars League"
$ws.Cells.Item(41, 3).Value = "Qatar Stars League"
$ws.Cells.Item(40, 4).Value = "Qatar Stars League"
$ws.Cells.Item(41, 4).Value = "Qatar Stars League"
$ws.Cells.Item(40, 5).Value = 45228.5625
$ws.Cells.Item(41, 5).Value = 45228.5625
$ws.Cells.Item(40, 6).Value = "AlWakrah SC"
$ws.Cells.Item(41, 6).Value = "Qatar SC Doha"
$ws.Cells.Item(40, 7).Value = "Umm Salal"
$ws.Cells.Item(41, 7).Value = "Al Markhiya"
$ws.Cells.Item(40, 8).Value = 2
$ws.Cells.Item(41, 8).Value = 4
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(40, 10).Value = "H"
$ws.Cells.Item(41, 10).Value = "H"
$ws.Cells.Item(40, 11).Value = 1.65
$ws.Cells.Item(41, 11).Value = 1.727
$ws.Cells.Item(40, 12).Value = 3.75
$ws.Cells.Item(41, 12).Value = 3.75
$ws.Cells.Item(40, 13).Value = 4.5
$ws.Cells.Item(41, 13).Value = 4
$ws.Cells.Item(40, 14).Value = 1.75
$ws.Cells.Item(41, 14).Value = 1.75
$ws.Cells.Item(40, 15).Value = 3.6
$ws.Cells.Item(41, 15).Value = 3.75
$ws.Cells.Item(40, 16).Value = 4
$ws.Cells.Item(41, 16).Value = 3.8
$ws.Cells.Item(40, 17).Value = -0.75
$ws.Cells.Item(41, 17).Value = -0.5
$ws.Cells.Item(40, 18).Value = 1.95
$ws.Cells.Item(41, 18).Value = 1.75
$ws.Cells.Item(40, 19).Value = 1.85
$ws.Cells.Item(41, 19).Value = 1.95
$ws.Cells.Item(40, 20).Value = 3
$ws.Cells.Item(41, 20).Value = 3
$ws.Cells.Item(40, 21).Value = 1.9
$ws.Cells.Item(41, 21).Value = 1.95
$ws.Cells.Item(40, 22).Value = 1.9
$ws.Cells.Item(41, 22).Value = 1.85
$ws.Cells.Item(40, 23).Value = 0.75
$ws.Cells.Item(41, 23).Value = 0.75
$ws.Cells.Item(40, 24).Value = -1
$ws.Cells.Item(41, 24).Value = -1
$ws.Cells.Item(40, 25).Value = -1
$ws.Cells.Item(41, 25).Value = -1
$ws.Cells.Item(40, 26).Value = 0.475
$ws.Cells.Item(41, 26).Value = 0.75
$ws.Cells.Item(40, 27).Value = -0.5
$ws.Cells.Item(41, 27).Value = -1
$ws.Cells.Item(40, 28).Value = 0
$ws.Cells.Item(41, 28).Value = 0.95
$ws.Cells.Item(40, 29).Value = -0
$ws.Cells.Item(41, 29).Value = -1

# Swap row 42 and row 43 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(42, 2).Value = 7003590
$ws.Cells.Item(43, 2).Value = 7004611
$ws.Cells.Item(42, 3).Value = "Qatar Stars League"
$ws.Cells.Item(43, 3).Value = "Qatar Stars League"
$ws.Cells.Item(42, 4).Value = "Qatar Stars League"
$ws.Cells.Item(43, 4).Value = "Qatar Stars League"
$ws.Cells.Item(42, 5).Value = 45232.47916666666
$ws.Cells.Item(43, 5).Value = 45232.47916666666
$ws.Cells.Item(42, 6).Value = "Al Sadd"
$ws.Cells.Item(43, 6).Value = "AlMuaidar"
$ws.Cells.Item(42, 7).Value = "Al Markhiya"
$ws.Cells.Item(43, 7).Value = "Umm Salal"
$ws.Cells.Item(42, 8).Value = 5
$ws.Cells.Item(43, 8).Value = 1
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(43, 9).Value = 3
$ws.Cells.Item(42, 10).Value = "H"
$ws.Cells.Item(43, 10).Value = "A"
$ws.Cells.Item(42, 11).Value = 1.125
$ws.Cells.Item(43, 11).Value = 4
$ws.Cells.Item(42, 12).Value = 8
$ws.Cells.Item(43, 12).Value = 3.25
$ws.Cells.Item(42, 13).Value = 10
$ws.Cells.Item(43, 13).Value = 1.8
$ws.Cells.Item(42, 14).Value = 1.2
$ws.Cells.Item(43, 14).Value = 3.75
$ws.Cells.Item(42, 15).Value = 5.75
$ws.Cells.Item(43, 15).Value = 3.4
$ws.Cells.Item(42, 16).Value = 10
$ws.Cells.Item(43, 16).Value = 1.8
$ws.Cells.Item(42, 17).Value = -2
$ws.Cells.Item(43, 17).Value = 0.5
$ws.Cells.Item(42, 18).Value = 1.975
$ws.Cells.Item(43, 18).Value = 1.925
$ws.Cells.Item(42, 19).Value = 1.825
$ws.Cells.Item(43, 19).Value = 1.875
$ws.Cells.Item(42, 20).Value = 3.5
$ws.Cells.Item(43, 20).Value = 3
$ws.Cells.Item(42, 21).Value = 1.95
$ws.Cells.Item(43, 21).Value = 1.85
$ws.Cells.Item(42, 22).Value = 1.85
$ws.Cells.Item(43, 22).Value = 1.95
$ws.Cells.Item(42, 23).Value = 0.2
$ws.Cells.Item(43, 23).Value = -1
$ws.Cells.Item(42, 24).Value = -1
$ws.Cells.Item(43, 24).Value = -1
$ws.Cells.Item(42, 25).Value = -1
$ws.Cells.Item(43, 25).Value = 0.8
$ws.Cells.Item(42, 26).Value = 0.9750000000000001
$ws.Cells.Item(43, 26).Value = -1
$ws.Cells.Item(42, 27).Value = -1
$ws.Cells.Item(43, 27).Value = 0.875
$ws.Cells.Item(42, 28).Value = 0.95
$ws.Cells.Item(43, 28).Value = 0.8500000000000001
$ws.Cells.Item(42, 29).Value = -1
$ws.Cells.Item(43, 29).Value = -1

# Swap row 56 and row 57 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(56, 2).Value = 7004619
$ws.Cells.Item(57, 2).Value = 7004618
$ws.Cells.Item(56, 3).Value = "Qatar Stars League"
$ws.Cells.Item(57, 3).Value = "Qatar Stars League"
$ws.Cells.Item(56, 4).Value = "Qatar Stars League"
$ws.Cells.Item(57, 4).Value = "Qatar Stars League"
$ws.Cells.Item(56, 5).Value = 45262.47916666666
$ws.Cells.Item(57, 5).Value = 45262.47916666666
$ws.Cells.Item(56, 6).Value = "Qatar SC Doha"
$ws.Cells.Item(57, 6).Value = "Al Gharafa"
$ws.Cells.Item(56, 7).Value = "AlMuaidar"
$ws.Cells.Item(57, 7).Value = "AlWakrah SC"
$ws.Cells.Item(56, 8).Value = 3
$ws.Cells.Item(57, 8).Value = 1
$ws.Cells.Item(56, 9).Value = 2
$ws.Cells.Item(57, 9).Value = 1
$ws.Cells.Item(56, 10).Value = "H"
$ws.Cells.Item(57, 10).Value = "D"
$ws.Cells.Item(56, 11).Value = 1.4
$ws.Cells.Item(57, 11).Value = 2.5
$ws.Cells.Item(56, 12).Value = 4.5
$ws.Cells.Item(57, 12).Value = 3.75
$ws.Cells.Item(56, 13).Value = 5.5
$ws.Cells.Item(57, 13).Value = 2.4
$ws.Cells.Item(56, 14).Value = 1.833
$ws.Cells.Item(57, 14).Value = 1.95
$ws.Cells.Item(56, 15).Value = 3.8
$ws.Cells.Item(57, 15).Value = 3.8
$ws.Cells.Item(56, 16).Value = 3.25
$ws.Cells.Item(57, 16).Value = 3.25
$ws.Cells.Item(56, 17).Value = -0.5
$ws.Cells.Item(57, 17).Value = -0.5
$ws.Cells.Item(56, 18).Value = 1.9
$ws.Cells.Item(57, 18).Value = 1.975
$ws.Cells.Item(56, 19).Value = 1.9
$ws.Cells.Item(57, 19).Value = 1.825
$ws.Cells.Item(56, 20).Value = 3
$ws.Cells.Item(57, 20).Value = 3.5
$ws.Cells.Item(56, 21).Value = 2
$ws.Cells.Item(57, 21).Value = 1.975
$ws.Cells.Item(56, 22).Value = 1.8
$ws.Cells.Item(57, 22).Value = 1.825
$ws.Cells.Item(56, 23).Value = 0.833
$ws.Cells.Item(57, 23).Value = -1
$ws.Cells.Item(56, 24).Value = -1
$ws.Cells.Item(57, 24).Value = 2.8
$ws.Cells.Item(56, 25).Value = -1
$ws.Cells.Item(57, 25).Value = -1
$ws.Cells.Item(56, 26).Value = 0.8999999999999999
$ws.Cells.Item(57, 26).Value = -1
$ws.Cells.Item(56, 27).Value = -1
$ws.Cells.Item(57, 27).Value = 0.825
$ws.Cells.Item(56, 28).Value = 1
$ws.Cells.Item(57, 28).Value = -1
$ws.Cells.Item(56, 29).Value = -1
$ws.Cells.Item(57, 29).Value = 0.825

# Swap row 62 and row 63 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(62, 2).Value = 7004622
$ws.Cells.Item(63, 2).Value = 7004621
$ws.Cells.Item(62, 3).Value = "Qatar Stars League"
$ws.Cells.Item(63, 3).Value = "Qatar Stars League"
$ws.Cells.Item(62, 4).Value = "Qatar Stars League"
$ws.Cells.Item(63, 4).Value = "Qatar Stars League"
$ws.Cells.Item(62, 5).Value = 45269.47916666666
$ws.Cells.Item(63, 5).Value = 45269.47916666666
$ws.Cells.Item(62, 6).Value = "Al Duhail"
$ws.Cells.Item(63, 6).Value = "Al Gharafa"
$ws.Cells.Item(62, 7).Value = "Qatar SC Doha"
$ws.Cells.Item(63, 7).Value = "AlShamal SC"
$ws.Cells.Item(62, 8).Value = 1
$ws.Cells.Item(63, 8).Value = 1
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(63, 9).Value = 1
$ws.Cells.Item(62, 10).Value = "D"
$ws.Cells.Item(63, 10).Value = "D"
$ws.Cells.Item(62, 11).Value = 1.5
$ws.Cells.Item(63, 11).Value = 1.5
$ws.Cells.Item(62, 12).Value = 4.5
$ws.Cells.Item(63, 12).Value = 4.5
$ws.Cells.Item(62, 13).Value = 5
$ws.Cells.Item(63, 13).Value = 5
$ws.Cells.Item(62, 14).Value = 1.8
$ws.Cells.Item(63, 14).Value = 1.363
$ws.Cells.Item(62, 15).Value = 4
$ws.Cells.Item(63, 15).Value = 5
$ws.Cells.Item(62, 16).Value = 3.6
$ws.Cells.Item(63, 16).Value = 6
$ws.Cells.Item(62, 17).Value = -0.5
$ws.Cells.Item(63, 17).Value = -1.5
$ws.Cells.Item(62, 18).Value = 1.8
$ws.Cells.Item(63, 18).Value = 1.975
$ws.Cells.Item(62, 19).Value = 2
$ws.Cells.Item(63, 19).Value = 1.825
$ws.Cells.Item(62, 20).Value = 3
$ws.Cells.Item(63, 20).Value = 3.75
$ws.Cells.Item(62, 21).Value = 1.8
$ws.Cells.Item(63, 21).Value = 1.975
$ws.Cells.Item(62, 22).Value = 2
$ws.Cells.Item(63, 22).Value = 1.825
$ws.Cells.Item(62, 23).Value = -1
$ws.Cells.Item(63, 23).Value = -1
$ws.Cells.Item(62, 24).Value = 3
$ws.Cells.Item(63, 24).Value = 4
$ws.Cells.Item(62, 25).Value = -1
$ws.Cells.Item(63, 25).Value = -1
$ws.Cells.Item(62, 26).Value = -1
$ws.Cells.Item(63, 26).Value = -1
$ws.Cells.Item(62, 27).Value = 1
$ws.Cells.Item(63, 27).Value = 0.825
$ws.Cells.Item(62, 28).Value = -1
$ws.Cells.Item(63, 28).Value = -1
$ws.Cells.Item(62, 29).Value = 1
$ws.Cells.Item(63, 29).Value = 0.825

# Swap row 68 and row 69 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(68, 2).Value = 7609336
$ws.Cells.Item(69, 2).Value = 7004626
$ws.Cells.Item(68, 3).Value = "Qatar Stars League"
$ws.Cells.Item(69, 3).Value = "Qatar Stars League"
$ws.Cells.Item(68, 4).Value = "Qatar Stars League"
$ws.Cells.Item(69, 4).Value = "Qatar Stars League"
$ws.Cells.Item(68, 5).Value = 45280.47916666666
$ws.Cells.Item(69, 5).Value = 45280.47916666666
$ws.Cells.Item(68, 6).Value = "AlMuaidar"
$ws.Cells.Item(69, 6).Value = "Al Gharafa"
$ws.Cells.Item(68, 7).Value = "AlWakrah SC"
$ws.Cells.Item(69, 7).Value = "Qatar SC Doha"
$ws.Cells.Item(68, 8).Value = 2
$ws.Cells.Item(69, 8).Value = 2
$ws.Cells.Item(68, 9).Value = 4
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(68, 10).Value = "A"
$ws.Cells.Item(69, 10).Value = "H"
$ws.Cells.Item(68, 11).Value = 4
$ws.Cells.Item(69, 11).Value = 1.909
$ws.Cells.Item(68, 12).Value = 4
$ws.Cells.Item(69, 12).Value = 3.8
$ws.Cells.Item(68, 13).Value = 1.65
$ws.Cells.Item(69, 13).Value = 3.25
$ws.Cells.Item(68, 14).Value = 4
$ws.Cells.Item(69, 14).Value = 1.909
$ws.Cells.Item(68, 15).Value = 3.75
$ws.Cells.Item(69, 15).Value = 3.75
$ws.Cells.Item(68, 16).Value = 1.666
$ws.Cells.Item(69, 16).Value = 3.4
$ws.Cells.Item(68, 17).Value = 0.75
$ws.Cells.Item(69, 17).Value = -0.5
$ws.Cells.Item(68, 18).Value = 1.925
$ws.Cells.Item(69, 18).Value = 1.925
$ws.Cells.Item(68, 19).Value = 1.875
$ws.Cells.Item(69, 19).Value = 1.875
$ws.Cells.Item(68, 20).Value = 3
$ws.Cells.Item(69, 20).Value = 3
$ws.Cells.Item(68, 21).Value = 1.825
$ws.Cells.Item(69, 21).Value = 1.825
$ws.Cells.Item(68, 22).Value = 1.975
$ws.Cells.Item(69, 22).Value = 1.975
$ws.Cells.Item(68, 23).Value = -1
$ws.Cells.Item(69, 23).Value = 0.909
$ws.Cells.Item(68, 24).Value = -1
$ws.Cells.Item(69, 24).Value = -1
$ws.Cells.Item(68, 25).Value = 0.6659999999999999
$ws.Cells.Item(69, 25).Value = -1
$ws.Cells.Item(68, 26).Value = -1
$ws.Cells.Item(69, 26).Value = 0.925
$ws.Cells.Item(68, 27).Value = 0.875
$ws.Cells.Item(69, 27).Value = -1
$ws.Cells.Item(68, 28).Value = 0.825
$ws.Cells.Item(69, 28).Value = 0
$ws.Cells.Item(68, 29).Value = -1
$ws.Cells.Item(69, 29).Value = -0

# Swap row 70 and row 71 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(70, 2).Value = 7609335
$ws.Cells.Item(71, 2).Value = 7004627
$ws.Cells.Item(70, 3).Value = "Qatar Stars League"
$ws.Cells.Item(71, 3).Value = "Qatar Stars League"
$ws.Cells.Item(70, 4).Value = "Qatar Stars League"
$ws.Cells.Item(71, 4).Value = "Qatar Stars League"
$ws.Cells.Item(70, 5).Value = 45280.5625
$ws.Cells.Item(71, 5).Value = 45280.5625
$ws.Cells.Item(70, 6).Value = "AlShamal SC"
$ws.Cells.Item(71, 6).Value = "AlRayyan SC"
$ws.Cells.Item(70, 7).Value = "AlArabi Doha"
$ws.Cells.Item(71, 7).Value = "Al Markhiya"
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(71, 8).Value = 6
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(70, 10).Value = "D"
$ws.Cells.Item(71, 10).Value = "H"
$ws.Cells.Item(70, 11).Value = 4.5
$ws.Cells.Item(71, 11).Value = 1.444
$ws.Cells.Item(70, 12).Value = 4.2
$ws.Cells.Item(71, 12).Value = 4.75
$ws.Cells.Item(70, 13).Value = 1.533
$ws.Cells.Item(71, 13).Value = 5.75
$ws.Cells.Item(70, 14).Value = 5.25
$ws.Cells.Item(71, 14).Value = 1.363
$ws.Cells.Item(70, 15).Value = 4.2
$ws.Cells.Item(71, 15).Value = 5
$ws.Cells.Item(70, 16).Value = 1.45
$ws.Cells.Item(71, 16).Value = 6.5
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(71, 17).Value = -1.5
$ws.Cells.Item(70, 18).Value = 2
$ws.Cells.Item(71, 18).Value = 2
$ws.Cells.Item(70, 19).Value = 1.8
$ws.Cells.Item(71, 19).Value = 1.8
$ws.Cells.Item(70, 20).Value = 3
$ws.Cells.Item(71, 20).Value = 3.25
$ws.Cells.Item(70, 21).Value = 1.9
$ws.Cells.Item(71, 21).Value = 1.95
$ws.Cells.Item(70, 22).Value = 1.9
$ws.Cells.Item(71, 22).Value = 1.85
$ws.Cells.Item(70, 23).Value = -1
$ws.Cells.Item(71, 23).Value = 0.363
$ws.Cells.Item(70, 24).Value = 3.2
$ws.Cells.Item(71, 24).Value = -1
$ws.Cells.Item(70, 25).Value = -1
$ws.Cells.Item(71, 25).Value = -1
$ws.Cells.Item(70, 26).Value = 1
$ws.Cells.Item(71, 26).Value = 1
$ws.Cells.Item(70, 27).Value = -1
$ws.Cells.Item(71, 27).Value = -1
$ws.Cells.Item(70, 28).Value = -1
$ws.Cells.Item(71, 28).Value = 0.95
$ws.Cells.Item(70, 29).Value = 0.8999999999999999
$ws.Cells.Item(71, 29).Value = -1

# Swap row 75 and row 76 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(75, 2).Value = 7840799
$ws.Cells.Item(76, 2).Value = 7840798
$ws.Cells.Item(75, 3).Value = "Qatar Stars League"
$ws.Cells.Item(76, 3).Value = "Qatar Stars League"
$ws.Cells.Item(75, 4).Value = "Qatar Stars League"
$ws.Cells.Item(76, 4).Value = "Qatar Stars League"
$ws.Cells.Item(75, 5).Value = 45346.58333333334
$ws.Cells.Item(76, 5).Value = 45346.58333333334
$ws.Cells.Item(75, 6).Value = "Al Markhiya"
$ws.Cells.Item(76, 6).Value = "Al Gharafa"
$ws.Cells.Item(75, 7).Value = "AlMuaidar"
$ws.Cells.Item(76, 7).Value = "AlAhli Doha"
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(76, 8).Value = 1
$ws.Cells.Item(75, 9).Value = 2
$ws.Cells.Item(76, 9).Value = 2
$ws.Cells.Item(75, 10).Value = "A"
$ws.Cells.Item(76, 10).Value = "A"
$ws.Cells.Item(75, 11).Value = 2.55
$ws.Cells.Item(76, 11).Value = 1.6
$ws.Cells.Item(75, 12).Value = 3.4
$ws.Cells.Item(76, 12).Value = 4.333
$ws.Cells.Item(75, 13).Value = 2.5
$ws.Cells.Item(76, 13).Value = 4.2
$ws.Cells.Item(75, 14).Value = 3.1
$ws.Cells.Item(76, 14).Value = 2.05
$ws.Cells.Item(75, 15).Value = 3.5
$ws.Cells.Item(76, 15).Value = 3.8
$ws.Cells.Item(75, 16).Value = 2.05
$ws.Cells.Item(76, 16).Value = 2.875
$ws.Cells.Item(75, 17).Value = 0.25
$ws.Cells.Item(76, 17).Value = -0.25
$ws.Cells.Item(75, 18).Value = 1.925
$ws.Cells.Item(76, 18).Value = 1.825
$ws.Cells.Item(75, 19).Value = 1.875
$ws.Cells.Item(76, 19).Value = 1.975
$ws.Cells.Item(75, 20).Value = 3
$ws.Cells.Item(76, 20).Value = 3.5
$ws.Cells.Item(75, 21).Value = 1.95
$ws.Cells.Item(76, 21).Value = 1.85
$ws.Cells.Item(75, 22).Value = 1.85
$ws.Cells.Item(76, 22).Value = 1.95
$ws.Cells.Item(75, 23).Value = -1
$ws.Cells.Item(76, 23).Value = -1
$ws.Cells.Item(75, 24).Value = -1
$ws.Cells.Item(76, 24).Value = -1
$ws.Cells.Item(75, 25).Value = 1.05
$ws.Cells.Item(76, 25).Value = 1.875
$ws.Cells.Item(75, 26).Value = -1
$ws.Cells.Item(76, 26).Value = -1
$ws.Cells.Item(75, 27).Value = 0.875
$ws.Cells.Item(76, 27).Value = 0.9750000000000001
$ws.Cells.Item(75, 28).Value = -1
$ws.Cells.Item(76, 28).Value = -1
$ws.Cells.Item(75, 29).Value = 0.8500000000000001
$ws.Cells.Item(76, 29).Value = 0.95

# Swap row 81 and row 82 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(81, 2).Value = 7840802
$ws.Cells.Item(82, 2).Value = 7840805
$ws.Cells.Item(81, 3).Value = "Qatar Stars League"
$ws.Cells.Item(82, 3).Value = "Qatar Stars League"
$ws.Cells.Item(81, 4).Value = "Qatar Stars League"
$ws.Cells.Item(82, 4).Value = "Qatar Stars League"
$ws.Cells.Item(81, 5).Value = 45351.58333333334
$ws.Cells.Item(82, 5).Value = 45351.58333333334
$ws.Cells.Item(81, 6).Value = "AlAhli Doha"
$ws.Cells.Item(82, 6).Value = "AlWakrah SC"
$ws.Cells.Item(81, 7).Value = "Umm Salal"
$ws.Cells.Item(82, 7).Value = "Al Markhiya"
$ws.Cells.Item(81, 8).Value = 1
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(81, 9).Value = 2
$ws.Cells.Item(82, 9).Value = 2
$ws.Cells.Item(81, 10).Value = "A"
$ws.Cells.Item(82, 10).Value = "A"
$ws.Cells.Item(81, 11).Value = 2.4
$ws.Cells.Item(82, 11).Value = 1.062
$ws.Cells.Item(81, 12).Value = 4
$ws.Cells.Item(82, 12).Value = 11
$ws.Cells.Item(81, 13).Value = 2.25
$ws.Cells.Item(82, 13).Value = 17
$ws.Cells.Item(81, 14).Value = 2.3
$ws.Cells.Item(82, 14).Value = 1.363
$ws.Cells.Item(81, 15).Value = 4
$ws.Cells.Item(82, 15).Value = 4.75
$ws.Cells.Item(81, 16).Value = 2.375
$ws.Cells.Item(82, 16).Value = 7
$ws.Cells.Item(81, 17).Value = 0
$ws.Cells.Item(82, 17).Value = -1.25
$ws.Cells.Item(81, 18).Value = 1.875
$ws.Cells.Item(82, 18).Value = 1.75
$ws.Cells.Item(81, 19).Value = 1.925
$ws.Cells.Item(82, 19).Value = 1.95
$ws.Cells.Item(81, 20).Value = 3.25
$ws.Cells.Item(82, 20).Value = 3.25
$ws.Cells.Item(81, 21).Value = 2
$ws.Cells.Item(82, 21).Value = 1.975
$ws.Cells.Item(81, 22).Value = 1.8
$ws.Cells.Item(82, 22).Value = 1.825
$ws.Cells.Item(81, 23).Value = -1
$ws.Cells.Item(82, 23).Value = -1
$ws.Cells.Item(81, 24).Value = -1
$ws.Cells.Item(82, 24).Value = -1
$ws.Cells.Item(81, 25).Value = 1.375
$ws.Cells.Item(82, 25).Value = 6
$ws.Cells.Item(81, 26).Value = -1
$ws.Cells.Item(82, 26).Value = -1
$ws.Cells.Item(81, 27).Value = 0.925
$ws.Cells.Item(82, 27).Value = 0.95
$ws.Cells.Item(81, 28).Value = -0.5
$ws.Cells.Item(82, 28).Value = -0.5
$ws.Cells.Item(81, 29).Value = 0.4
$ws.Cells.Item(82, 29).Value = 0.4125

# Swap row 83 and row 84 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(83, 2).Value = 7840803
$ws.Cells.Item(84, 2).Value = 7840685
$ws.Cells.Item(83, 3).Value = "Qatar Stars League"
$ws.Cells.Item(84, 3).Value = "Qatar Stars League"
$ws.Cells.Item(83, 4).Value = "Qatar Stars League"
$ws.Cells.Item(84, 4).Value = "Qatar Stars League"
$ws.Cells.Item(83, 5).Value = 45352.5
$ws.Cells.Item(84, 5).Value = 45352.5
$ws.Cells.Item(83, 6).Value = "Al Sadd"
$ws.Cells.Item(84, 6).Value = "Al Duhail"
$ws.Cells.Item(83, 7).Value = "Qatar SC Doha"
$ws.Cells.Item(84, 7).Value = "Al Gharafa"
$ws.Cells.Item(83, 8).Value = 3
$ws.Cells.Item(84, 8).Value = 1
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(84, 9).Value = 4
$ws.Cells.Item(83, 10).Value = "H"
$ws.Cells.Item(84, 10).Value = "A"
$ws.Cells.Item(83, 11).Value = 1.285
$ws.Cells.Item(84, 11).Value = 2.1
$ws.Cells.Item(83, 12).Value = 5.75
$ws.Cells.Item(84, 12).Value = 3.8
$ws.Cells.Item(83, 13).Value = 7.5
$ws.Cells.Item(84, 13).Value = 2.9
$ws.Cells.Item(83, 14).Value = 1.2
$ws.Cells.Item(84, 14).Value = 1.8
$ws.Cells.Item(83, 15).Value = 6.5
$ws.Cells.Item(84, 15).Value = 4
$ws.Cells.Item(83, 16).Value = 9.5
$ws.Cells.Item(84, 16).Value = 3.6
$ws.Cells.Item(83, 17).Value = -2
$ws.Cells.Item(84, 17).Value = -0.75
$ws.Cells.Item(83, 18).Value = 1.95
$ws.Cells.Item(84, 18).Value = 1.975
$ws.Cells.Item(83, 19).Value = 1.85
$ws.Cells.Item(84, 19).Value = 1.825
$ws.Cells.Item(83, 20).Value = 3.5
$ws.Cells.Item(84, 20).Value = 3.5
$ws.Cells.Item(83, 21).Value = 1.95
$ws.Cells.Item(84, 21).Value = 1.825
$ws.Cells.Item(83, 22).Value = 1.85
$ws.Cells.Item(84, 22).Value = 1.975
$ws.Cells.Item(83, 23).Value = 0.2
$ws.Cells.Item(84, 23).Value = -1
$ws.Cells.Item(83, 24).Value = -1
$ws.Cells.Item(84, 24).Value = -1
$ws.Cells.Item(83, 25).Value = -1
$ws.Cells.Item(84, 25).Value = 2.6
$ws.Cells.Item(83, 26).Value = 0.95
$ws.Cells.Item(84, 26).Value = -1
$ws.Cells.Item(83, 27).Value = -1
$ws.Cells.Item(84, 27).Value = 0.825
$ws.Cells.Item(83, 28).Value = -1
$ws.Cells.Item(84, 28).Value = 0.825
$ws.Cells.Item(83, 29).Value = 0.8500000000000001
$ws.Cells.Item(84, 29).Value = -1

# Swap row 90 and row 91 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(90, 2).Value = 7840810
$ws.Cells.Item(91, 2).Value = 7840809
$ws.Cells.Item(90, 3).Value = "Qatar Stars League"
$ws.Cells.Item(91, 3).Value = "Qatar Stars League"
$ws.Cells.Item(90, 4).Value = "Qatar Stars League"
$ws.Cells.Item(91, 4).Value = "Qatar Stars League"
$ws.Cells.Item(90, 5).Value = 45357.58333333334
$ws.Cells.Item(91, 5).Value = 45357.58333333334
$ws.Cells.Item(90, 6).Value = "AlRayyan SC"
$ws.Cells.Item(91, 6).Value = "Al Gharafa"
$ws.Cells.Item(90, 7).Value = "AlWakrah SC"
$ws.Cells.Item(91, 7).Value = "Al Sadd"
$ws.Cells.Item(90, 8).Value = 3
$ws.Cells.Item(91, 8).Value = 2
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(91, 9).Value = 2
$ws.Cells.Item(90, 10).Value = "H"
$ws.Cells.Item(91, 10).Value = "D"
$ws.Cells.Item(90, 11).Value = 2
$ws.Cells.Item(91, 11).Value = 5
$ws.Cells.Item(90, 12).Value = 3.6
$ws.Cells.Item(91, 12).Value = 4.75
$ws.Cells.Item(90, 13).Value = 3.1
$ws.Cells.Item(91, 13).Value = 1.45
$ws.Cells.Item(90, 14).Value = 2.15
$ws.Cells.Item(91, 14).Value = 5.25
$ws.Cells.Item(90, 15).Value = 3.4
$ws.Cells.Item(91, 15).Value = 5
$ws.Cells.Item(90, 16).Value = 2.9
$ws.Cells.Item(91, 16).Value = 1.4
$ws.Cells.Item(90, 17).Value = -0.25
$ws.Cells.Item(91, 17).Value = 1.25
$ws.Cells.Item(90, 18).Value = 1.975
$ws.Cells.Item(91, 18).Value = 2
$ws.Cells.Item(90, 19).Value = 1.825
$ws.Cells.Item(91, 19).Value = 1.8
$ws.Cells.Item(90, 20).Value = 3
$ws.Cells.Item(91, 20).Value = 3.75
$ws.Cells.Item(90, 21).Value = 1.925
$ws.Cells.Item(91, 21).Value = 1.875
$ws.Cells.Item(90, 22).Value = 1.875
$ws.Cells.Item(91, 22).Value = 1.925
$ws.Cells.Item(90, 23).Value = 1.15
$ws.Cells.Item(91, 23).Value = -1
$ws.Cells.Item(90, 24).Value = -1
$ws.Cells.Item(91, 24).Value = 4
$ws.Cells.Item(90, 25).Value = -1
$ws.Cells.Item(91, 25).Value = -1
$ws.Cells.Item(90, 26).Value = 0.9750000000000001
$ws.Cells.Item(91, 26).Value = 1
$ws.Cells.Item(90, 27).Value = -1
$ws.Cells.Item(91, 27).Value = -1
$ws.Cells.Item(90, 28).Value = 0
$ws.Cells.Item(91, 28).Value = 0.4375
$ws.Cells.Item(90, 29).Value = -0
$ws.Cells.Item(91, 29).Value = -0.5

# Swap row 98 and row 99 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(98, 2).Value = 7840688
$ws.Cells.Item(99, 2).Value = 7840816
$ws.Cells.Item(98, 3).Value = "Qatar Stars League"
$ws.Cells.Item(99, 3).Value = "Qatar Stars League"
$ws.Cells.Item(98, 4).Value = "Qatar Stars League"
$ws.Cells.Item(99, 4).Value = "Qatar Stars League"
$ws.Cells.Item(98, 5).Value = 45366.64583333334
$ws.Cells.Item(99, 5).Value = 45366.64583333334
$ws.Cells.Item(98, 6).Value = "Al Duhail"
$ws.Cells.Item(99, 6).Value = "Al Sadd"
$ws.Cells.Item(98, 7).Value = "AlShamal SC"
$ws.Cells.Item(99, 7).Value = "AlMuaidar"
$ws.Cells.Item(98, 8).Value = 3
$ws.Cells.Item(99, 8).Value = 4
$ws.Cells.Item(98, 9).Value = 1
$ws.Cells.Item(99, 9).Value = 2
$ws.Cells.Item(98, 10).Value = "H"
$ws.Cells.Item(99, 10).Value = "H"
$ws.Cells.Item(98, 11).Value = 1.444
$ws.Cells.Item(99, 11).Value = 1.166
$ws.Cells.Item(98, 12).Value = 4.75
$ws.Cells.Item(99, 12).Value = 7
$ws.Cells.Item(98, 13).Value = 5.5
$ws.Cells.Item(99, 13).Value = 11
$ws.Cells.Item(98, 14).Value = 1.615
$ws.Cells.Item(99, 14).Value = 1.285
$ws.Cells.Item(98, 15).Value = 4.2
$ws.Cells.Item(99, 15).Value = 5.5
$ws.Cells.Item(98, 16).Value = 4.333
$ws.Cells.Item(99, 16).Value = 8
$ws.Cells.Item(98, 17).Value = -0.75
$ws.Cells.Item(99, 17).Value = -1.75
$ws.Cells.Item(98, 18).Value = 1.8
$ws.Cells.Item(99, 18).Value = 1.975
$ws.Cells.Item(98, 19).Value = 2
$ws.Cells.Item(99, 19).Value = 1.825
$ws.Cells.Item(98, 20).Value = 3.25
$ws.Cells.Item(99, 20).Value = 3.75
$ws.Cells.Item(98, 21).Value = 1.925
$ws.Cells.Item(99, 21).Value = 1.925
$ws.Cells.Item(98, 22).Value = 1.875
$ws.Cells.Item(99, 22).Value = 1.775
$ws.Cells.Item(98, 23).Value = 0.615
$ws.Cells.Item(99, 23).Value = 0.2849999999999999
$ws.Cells.Item(98, 24).Value = -1
$ws.Cells.Item(99, 24).Value = -1
$ws.Cells.Item(98, 25).Value = -1
$ws.Cells.Item(99, 25).Value = -1
$ws.Cells.Item(98, 26).Value = 0.8
$ws.Cells.Item(99, 26).Value = 0.4875
$ws.Cells.Item(98, 27).Value = -1
$ws.Cells.Item(99, 27).Value = -0.5
$ws.Cells.Item(98, 28).Value = 0.925
$ws.Cells.Item(99, 28).Value = 0.925
$ws.Cells.Item(98, 29).Value = -1
$ws.Cells.Item(99, 29).Value = -1

# Swap row 104 and row 105 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(104, 2).Value = 7004654
$ws.Cells.Item(105, 2).Value = 7003493
$ws.Cells.Item(104, 3).Value = "Qatar Stars League"
$ws.Cells.Item(105, 3).Value = "Qatar Stars League"
$ws.Cells.Item(104, 4).Value = "Qatar Stars League"
$ws.Cells.Item(105, 4).Value = "Qatar Stars League"
$ws.Cells.Item(104, 5).Value = 45380.64583333334
$ws.Cells.Item(105, 5).Value = 45380.64583333334
$ws.Cells.Item(104, 6).Value = "Umm Salal"
$ws.Cells.Item(105, 6).Value = "AlShamal SC"
$ws.Cells.Item(104, 7).Value = "AlMuaidar"
$ws.Cells.Item(105, 7).Value = "AlAhli Doha"
$ws.Cells.Item(104, 8).Value = 1
$ws.Cells.Item(105, 8).Value = 2
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(104, 10).Value = "A"
$ws.Cells.Item(105, 10).Value = "H"
$ws.Cells.Item(104, 11).Value = 2.7
$ws.Cells.Item(105, 11).Value = 2.55
$ws.Cells.Item(104, 12).Value = 3.5
$ws.Cells.Item(105, 12).Value = 3.5
$ws.Cells.Item(104, 13).Value = 2.25
$ws.Cells.Item(105, 13).Value = 2.375
$ws.Cells.Item(104, 14).Value = 2.7
$ws.Cells.Item(105, 14).Value = 2.875
$ws.Cells.Item(104, 15).Value = 3.5
$ws.Cells.Item(105, 15).Value = 3.5
$ws.Cells.Item(104, 16).Value = 2.25
$ws.Cells.Item(105, 16).Value = 2.15
$ws.Cells.Item(104, 17).Value = 0.25
$ws.Cells.Item(105, 17).Value = 0.25
$ws.Cells.Item(104, 18).Value = 1.825
$ws.Cells.Item(105, 18).Value = 1.825
$ws.Cells.Item(104, 19).Value = 1.975
$ws.Cells.Item(105, 19).Value = 1.975
$ws.Cells.Item(104, 20).Value = 3
$ws.Cells.Item(105, 20).Value = 3
$ws.Cells.Item(104, 21).Value = 1.925
$ws.Cells.Item(105, 21).Value = 1.875
$ws.Cells.Item(104, 22).Value = 1.875
$ws.Cells.Item(105, 22).Value = 1.925
$ws.Cells.Item(104, 23).Value = -1
$ws.Cells.Item(105, 23).Value = 1.875
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(104, 25).Value = 1.25
$ws.Cells.Item(105, 25).Value = -1
$ws.Cells.Item(104, 26).Value = -1
$ws.Cells.Item(105, 26).Value = 0.825
$ws.Cells.Item(104, 27).Value = 0.9750000000000001
$ws.Cells.Item(105, 27).Value = -1
$ws.Cells.Item(104, 28).Value = 0
$ws.Cells.Item(105, 28).Value = 0
$ws.Cells.Item(104, 29).Value = -0
$ws.Cells.Item(105, 29).Value = -0

# Swap row 106 and row 107 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(106, 2).Value = 7004653
$ws.Cells.Item(107, 2).Value = 7004656
$ws.Cells.Item(106, 3).Value = "Qatar Stars League"
$ws.Cells.Item(107, 3).Value = "Qatar Stars League"
$ws.Cells.Item(106, 4).Value = "Qatar Stars League"
$ws.Cells.Item(107, 4).Value = "Qatar Stars League"
$ws.Cells.Item(106, 5).Value = 45381.64583333334
$ws.Cells.Item(107, 5).Value = 45381.64583333334
$ws.Cells.Item(106, 6).Value = "Qatar SC Doha"
$ws.Cells.Item(107, 6).Value = "AlArabi Doha"
$ws.Cells.Item(106, 7).Value = "AlWakrah SC"
$ws.Cells.Item(107, 7).Value = "Al Gharafa"
$ws.Cells.Item(106, 8).Value = 1
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 5
$ws.Cells.Item(107, 9).Value = 1
$ws.Cells.Item(106, 10).Value = "A"
$ws.Cells.Item(107, 10).Value = "A"
$ws.Cells.Item(106, 11).Value = 3.5
$ws.Cells.Item(107, 11).Value = 2.55
$ws.Cells.Item(106, 12).Value = 3.6
$ws.Cells.Item(107, 12).Value = 3.6
$ws.Cells.Item(106, 13).Value = 1.85
$ws.Cells.Item(107, 13).Value = 2.3
$ws.Cells.Item(106, 14).Value = 3.4
$ws.Cells.Item(107, 14).Value = 2.4
$ws.Cells.Item(106, 15).Value = 3.6
$ws.Cells.Item(107, 15).Value = 3.6
$ws.Cells.Item(106, 16).Value = 1.85
$ws.Cells.Item(107, 16).Value = 2.375
$ws.Cells.Item(106, 17).Value = 0.5
$ws.Cells.Item(107, 17).Value = 0
$ws.Cells.Item(106, 18).Value = 1.9
$ws.Cells.Item(107, 18).Value = 1.95
$ws.Cells.Item(106, 19).Value = 1.9
$ws.Cells.Item(107, 19).Value = 1.85
$ws.Cells.Item(106, 20).Value = 3
$ws.Cells.Item(107, 20).Value = 3.25
$ws.Cells.Item(106, 21).Value = 1.925
$ws.Cells.Item(107, 21).Value = 1.775
$ws.Cells.Item(106, 22).Value = 1.875
$ws.Cells.Item(107, 22).Value = 1.925
$ws.Cells.Item(106, 23).Value = -1
$ws.Cells.Item(107, 23).Value = -1
$ws.Cells.Item(106, 24).Value = -1
$ws.Cells.Item(107, 24).Value = -1
$ws.Cells.Item(106, 25).Value = 0.8500000000000001
$ws.Cells.Item(107, 25).Value = 1.375
$ws.Cells.Item(106, 26).Value = -1
$ws.Cells.Item(107, 26).Value = -1
$ws.Cells.Item(106, 27).Value = 0.8999999999999999
$ws.Cells.Item(107, 27).Value = 0.8500000000000001
$ws.Cells.Item(106, 28).Value = 0.925
$ws.Cells.Item(107, 28).Value = -1
$ws.Cells.Item(106, 29).Value = -1
$ws.Cells.Item(107, 29).Value = 0.925

# Swap row 108 and row 109 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(108, 2).Value = 7882227
$ws.Cells.Item(109, 2).Value = 7004655
$ws.Cells.Item(108, 3).Value = "Qatar Stars League"
$ws.Cells.Item(109, 3).Value = "Qatar Stars League"
$ws.Cells.Item(108, 4).Value = "Qatar Stars League"
$ws.Cells.Item(109, 4).Value = "Qatar Stars League"
$ws.Cells.Item(108, 5).Value = 45382.64583333334
$ws.Cells.Item(109, 5).Value = 45382.64583333334
$ws.Cells.Item(108, 6).Value = "Al Markhiya"
$ws.Cells.Item(109, 6).Value = "AlRayyan SC"
$ws.Cells.Item(108, 7).Value = "Al Sadd"
$ws.Cells.Item(109, 7).Value = "Al Duhail"
$ws.Cells.Item(108, 8).Value = 1
$ws.Cells.Item(109, 8).Value = 2
$ws.Cells.Item(108, 9).Value = 2
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(108, 10).Value = "A"
$ws.Cells.Item(109, 10).Value = "H"
$ws.Cells.Item(108, 11).Value = 10
$ws.Cells.Item(109, 11).Value = 2.1
$ws.Cells.Item(108, 12).Value = 6.5
$ws.Cells.Item(109, 12).Value = 3.75
$ws.Cells.Item(108, 13).Value = 1.2
$ws.Cells.Item(109, 13).Value = 2.8
$ws.Cells.Item(108, 14).Value = 11
$ws.Cells.Item(109, 14).Value = 2.25
$ws.Cells.Item(108, 15).Value = 7.5
$ws.Cells.Item(109, 15).Value = 3.75
$ws.Cells.Item(108, 16).Value = 1.142
$ws.Cells.Item(109, 16).Value = 2.6
$ws.Cells.Item(108, 17).Value = 2.25
$ws.Cells.Item(109, 17).Value = -0.25
$ws.Cells.Item(108, 18).Value = 1.95
$ws.Cells.Item(109, 18).Value = 1.975
$ws.Cells.Item(108, 19).Value = 1.85
$ws.Cells.Item(109, 19).Value = 1.825
$ws.Cells.Item(108, 20).Value = 3.75
$ws.Cells.Item(109, 20).Value = 3.5
$ws.Cells.Item(108, 21).Value = 1.825
$ws.Cells.Item(109, 21).Value = 1.925
$ws.Cells.Item(108, 22).Value = 1.975
$ws.Cells.Item(109, 22).Value = 1.775
$ws.Cells.Item(108, 23).Value = -1
$ws.Cells.Item(109, 23).Value = 1.25
$ws.Cells.Item(108, 24).Value = -1
$ws.Cells.Item(109, 24).Value = -1
$ws.Cells.Item(108, 25).Value = 0.1419999999999999
$ws.Cells.Item(109, 25).Value = -1
$ws.Cells.Item(108, 26).Value = 0.95
$ws.Cells.Item(109, 26).Value = 0.9750000000000001
$ws.Cells.Item(108, 27).Value = -1
$ws.Cells.Item(109, 27).Value = -1
$ws.Cells.Item(108, 28).Value = -1
$ws.Cells.Item(109, 28).Value = -1
$ws.Cells.Item(108, 29).Value = 0.9750000000000001
$ws.Cells.Item(109, 29).Value = 0.7749999999999999

# Swap row 110 and row 111 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(110, 2).Value = 8022181
$ws.Cells.Item(111, 2).Value = 7004660
$ws.Cells.Item(110, 3).Value = "Qatar Stars League"
$ws.Cells.Item(111, 3).Value = "Qatar Stars League"
$ws.Cells.Item(110, 4).Value = "Qatar Stars League"
$ws.Cells.Item(111, 4).Value = "Qatar Stars League"
$ws.Cells.Item(110, 5).Value = 45387.64583333334
$ws.Cells.Item(111, 5).Value = 45387.64583333334
$ws.Cells.Item(110, 6).Value = "Al Duhail"
$ws.Cells.Item(111, 6).Value = "AlWakrah SC"
$ws.Cells.Item(110, 7).Value = "Al Sadd"
$ws.Cells.Item(111, 7).Value = "AlArabi Doha"
$ws.Cells.Item(110, 8).Value = 3
$ws.Cells.Item(111, 8).Value = 2
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(111, 9).Value = 4
$ws.Cells.Item(110, 10).Value = "H"
$ws.Cells.Item(111, 10).Value = "A"
$ws.Cells.Item(110, 11).Value = 4.5
$ws.Cells.Item(111, 11).Value = 2
$ws.Cells.Item(110, 12).Value = 4.333
$ws.Cells.Item(111, 12).Value = 3.75
$ws.Cells.Item(110, 13).Value = 1.6
$ws.Cells.Item(111, 13).Value = 3.2
$ws.Cells.Item(110, 14).Value = 4
$ws.Cells.Item(111, 14).Value = 1.909
$ws.Cells.Item(110, 15).Value = 4.2
$ws.Cells.Item(111, 15).Value = 3.8
$ws.Cells.Item(110, 16).Value = 1.666
$ws.Cells.Item(111, 16).Value = 3.4
$ws.Cells.Item(110, 17).Value = 0.75
$ws.Cells.Item(111, 17).Value = -0.5
$ws.Cells.Item(110, 18).Value = 1.95
$ws.Cells.Item(111, 18).Value = 1.95
$ws.Cells.Item(110, 19).Value = 1.85
$ws.Cells.Item(111, 19).Value = 1.85
$ws.Cells.Item(110, 20).Value = 3.5
$ws.Cells.Item(111, 20).Value = 3.25
$ws.Cells.Item(110, 21).Value = 1.85
$ws.Cells.Item(111, 21).Value = 2
$ws.Cells.Item(110, 22).Value = 1.95
$ws.Cells.Item(111, 22).Value = 1.8
$ws.Cells.Item(110, 23).Value = 3
$ws.Cells.Item(111, 23).Value = -1
$ws.Cells.Item(110, 24).Value = -1
$ws.Cells.Item(111, 24).Value = -1
$ws.Cells.Item(110, 25).Value = -1
$ws.Cells.Item(111, 25).Value = 2.4
$ws.Cells.Item(110, 26).Value = 0.95
$ws.Cells.Item(111, 26).Value = -1
$ws.Cells.Item(110, 27).Value = -1
$ws.Cells.Item(111, 27).Value = 0.8500000000000001
$ws.Cells.Item(110, 28).Value = 0.8500000000000001
$ws.Cells.Item(111, 28).Value = 1
$ws.Cells.Item(110, 29).Value = -1
$ws.Cells.Item(111, 29).Value = -1

# Swap row 118 and row 119 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(118, 2).Value = 7818846
$ws.Cells.Item(119, 2).Value = 7818294
$ws.Cells.Item(118, 3).Value = "Qatar Stars League"
$ws.Cells.Item(119, 3).Value = "Qatar Stars League"
$ws.Cells.Item(118, 4).Value = "Qatar Stars League"
$ws.Cells.Item(119, 4).Value = "Qatar Stars League"
$ws.Cells.Item(118, 5).Value = 45399.52083333334
$ws.Cells.Item(119, 5).Value = 45399.52083333334
$ws.Cells.Item(118, 6).Value = "AlMuaidar"
$ws.Cells.Item(119, 6).Value = "Umm Salal"
$ws.Cells.Item(118, 7).Value = "Al Duhail"
$ws.Cells.Item(119, 7).Value = "AlArabi Doha"
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(119, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 5
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(118, 10).Value = "A"
$ws.Cells.Item(119, 10).Value = "H"
$ws.Cells.Item(118, 11).Value = 4
$ws.Cells.Item(119, 11).Value = 3.75
$ws.Cells.Item(118, 12).Value = 4
$ws.Cells.Item(119, 12).Value = 3.6
$ws.Cells.Item(118, 13).Value = 1.666
$ws.Cells.Item(119, 13).Value = 1.8
$ws.Cells.Item(118, 14).Value = 4
$ws.Cells.Item(119, 14).Value = 4.333
$ws.Cells.Item(118, 15).Value = 3.8
$ws.Cells.Item(119, 15).Value = 3.8
$ws.Cells.Item(118, 16).Value = 1.7
$ws.Cells.Item(119, 16).Value = 1.65
$ws.Cells.Item(118, 17).Value = 0.75
$ws.Cells.Item(119, 17).Value = 0.75
$ws.Cells.Item(118, 18).Value = 1.9
$ws.Cells.Item(119, 18).Value = 1.95
$ws.Cells.Item(118, 19).Value = 1.9
$ws.Cells.Item(119, 19).Value = 1.85
$ws.Cells.Item(118, 20).Value = 3.25
$ws.Cells.Item(119, 20).Value = 3.25
$ws.Cells.Item(118, 21).Value = 2
$ws.Cells.Item(119, 21).Value = 1.975
$ws.Cells.Item(118, 22).Value = 1.8
$ws.Cells.Item(119, 22).Value = 1.825
$ws.Cells.Item(118, 23).Value = -1
$ws.Cells.Item(119, 23).Value = 3.333
$ws.Cells.Item(118, 24).Value = -1
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(118, 25).Value = 0.7
$ws.Cells.Item(119, 25).Value = -1
$ws.Cells.Item(118, 26).Value = -1
$ws.Cells.Item(119, 26).Value = 0.95
$ws.Cells.Item(118, 27).Value = 0.8999999999999999
$ws.Cells.Item(119, 27).Value = -1
$ws.Cells.Item(118, 28).Value = 1
$ws.Cells.Item(119, 28).Value = -0.5
$ws.Cells.Item(118, 29).Value = -1
$ws.Cells.Item(119, 29).Value = 0.4125

# Swap row 120 and row 121 (columns B:AC); column A (sequence id) stays fixed
$ws.Cells.Item(120, 2).Value = 7818282
$ws.Cells.Item(121, 2).Value = 7818593
$ws.Cells.Item(120, 3).Value = "Qatar Stars League"
$ws.Cells.Item(121, 3).Value = "Qatar Stars League"
$ws.Cells.Item(120, 4).Value = "Qatar Stars League"
$ws.Cells.Item(121, 4).Value = "Qatar Stars League"
$ws.Cells.Item(120, 5).Value = 45399.60416666666
$ws.Cells.Item(121, 5).Value = 45399.60416666666
$ws.Cells.Item(120, 6).Value = "Qatar SC Doha"
$ws.Cells.Item(121, 6).Value = "AlAhli Doha"
$ws.Cells.Item(120, 7).Value = "AlRayyan SC"
$ws.Cells.Item(121, 7).Value = "Al Sadd"
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(121, 8).Value = 1
$ws.Cells.Item(120, 9).Value = 3
$ws.Cells.Item(121, 9).Value = 9
$ws.Cells.Item(120, 10).Value = "A"
$ws.Cells.Item(121, 10).Value = "A"
$ws.Cells.Item(120, 11).Value = 3.6
$ws.Cells.Item(121, 11).Value = 8
$ws.Cells.Item(120, 12).Value = 3.6
$ws.Cells.Item(121, 12).Value = 6
$ws.Cells.Item(120, 13).Value = 1.8
$ws.Cells.Item(121, 13).Value = 1.222
$ws.Cells.Item(120, 14).Value = 4.333
$ws.Cells.Item(121, 14).Value = 9.5
$ws.Cells.Item(120, 15).Value = 4
$ws.Cells.Item(121, 15).Value = 7
$ws.Cells.Item(120, 16).Value = 1.6
$ws.Cells.Item(121, 16).Value = 1.166
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(121, 17).Value = 2.25
$ws.Cells.Item(120, 18).Value = 1.75
$ws.Cells.Item(121, 18).Value = 1.825
$ws.Cells.Item(120, 19).Value = 1.95
$ws.Cells.Item(121, 19).Value = 1.975
$ws.Cells.Item(120, 20).Value = 3.25
$ws.Cells.Item(121, 20).Value = 4
$ws.Cells.Item(120, 21).Value = 1.925
$ws.Cells.Item(121, 21).Value = 1.825
$ws.Cells.Item(120, 22).Value = 1.875
$ws.Cells.Item(121, 22).Value = 1.975
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(121, 23).Value = -1
$ws.Cells.Item(120, 24).Value = -1
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(120, 25).Value = 0.6000000000000001
$ws.Cells.Item(121, 25).Value = 0.1659999999999999
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.95
$ws.Cells.Item(121, 27).Value = 0.9750000000000001
$ws.Cells.Item(120, 28).Value = -0.5
$ws.Cells.Item(121, 28).Value = 0.825
$ws.Cells.Item(120, 29).Value = 0.4375
$ws.Cells.Item(121, 29).Value = -1
